$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 96.92308
$ws.Range("I33").Value = 103.63636
$ws.Range("K33").Value = 103.63636
$ws.Range("M33").Value = 125.36364

$ws.Range("H40").Value = 5756.294
$ws.Range("I40").Value = 4241.909
$ws.Range("J40").Value = 8532.666999999999
$ws.Range("K40").Value = 4241.909
$ws.Range("L40").Value = 8532.666999999999
$ws.Range("M40").Value = -4066.909
$ws.Range("N40").Value = -8882.666999999999

$ws.Range("H41").Value = 898.3333
$ws.Range("I41").Value = 175
$ws.Range("J41").Value = 2345
$ws.Range("K41").Value = 175
$ws.Range("L41").Value = 2345
$ws.Range("M41").Value = 265
$ws.Range("N41").Value = -3225

$ws.Range("H43").Value = 2000
$ws.Range("I43").Value = 2000
$ws.Range("K43").Value = 2000
$ws.Range("M43").Value = -1931

$ws.Range("H69").Value = 7180.8486

$ws.Range("H72").Value = 7180.8486

$ws.Range("H74").Value = 10107.143
$ws.Range("I74").Value = 4150
$ws.Range("K74").Value = 4150
$ws.Range("M74").Value = -3214

$ws.Range("H77").Value = 10107.143
$ws.Range("I77").Value = 4150
$ws.Range("K77").Value = 20750
$ws.Range("M77").Value = -16070

$ws.Range("H82").Value = 689.5
$ws.Range("I82").Value = 689.5
$ws.Range("K82").Value = 2068.5
$ws.Range("M82").Value = -1662.5

$ws.Range("H85").Value = 689.5
$ws.Range("I85").Value = 689.5
$ws.Range("K85").Value = 2068.5
$ws.Range("M85").Value = -664.5

$ws.Range("H107").Value = 107.26667
$ws.Range("I107").Value = 79.28570999999999
$ws.Range("K107").Value = 79.28570999999999
$ws.Range("M107").Value = 1840.71429

$ws.Range("H121").Value = 781.1667
$ws.Range("J121").Value = 781.1667
$ws.Range("L121").Value = 2343.5001
$ws.Range("N121").Value = -5837.5001

$ws.Range("H129").Value = 1734.3529
$ws.Range("I129").Value = 1153.6
$ws.Range("K129").Value = 3460.8
$ws.Range("M129").Value = 1539.2

$ws.Range("H131").Value = 1283.375
$ws.Range("I131").Value = 743.75
$ws.Range("K131").Value = 2231.25
$ws.Range("M131").Value = 2808.75

$ws.Range("H140").Value = 76195
$ws.Range("J140").Value = 94926.664
$ws.Range("L140").Value = 94926.664
$ws.Range("N140").Value = -105286.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7109.5264
$ws.Range("I32").Value = 6949
$ws.Range("K32").Value = 6949
$ws.Range("M32").Value = -6662

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H106").Value = 29184.5
$ws.Range("J106").Value = 29184.5
$ws.Range("L106").Value = 29184.5
$ws.Range("N106").Value = -31708.5

$ws.Range("H110").Value = 4722.857
$ws.Range("I110").Value = 4702.5
$ws.Range("J110").Value = 4750
$ws.Range("K110").Value = 4702.5
$ws.Range("L110").Value = 4750
$ws.Range("M110").Value = -2657.5
$ws.Range("N110").Value = -8840

$ws.Range("H122").Value = 1402.8462
$ws.Range("I122").Value = 1226.091
$ws.Range("K122").Value = 3678.273
$ws.Range("M122").Value = -1228.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 18498.4
$ws.Range("J88").Value = 18498.4
$ws.Range("L88").Value = 18498.4
$ws.Range("N88").Value = -19310.4

$ws.Range("H91").Value = 18498.4
$ws.Range("J91").Value = 18498.4
$ws.Range("L91").Value = 18498.4
$ws.Range("N91").Value = -21306.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2307.64
$ws.Range("I31").Value = 1677.9546
$ws.Range("K31").Value = 1677.9546
$ws.Range("M31").Value = -1382.9546

$ws.Range("H34").Value = 2307.64
$ws.Range("I34").Value = 1677.9546
$ws.Range("K34").Value = 1677.9546
$ws.Range("M34").Value = -1475.9546

$ws.Range("H59").Value = 26880.4
$ws.Range("J59").Value = 27644.889
$ws.Range("L59").Value = 27644.889
$ws.Range("N59").Value = -29934.889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 217.16667
$ws.Range("I86").Value = 209.4
$ws.Range("J86").Value = 256
$ws.Range("K86").Value = 628.2
$ws.Range("L86").Value = 768
$ws.Range("M86").Value = 557.8
$ws.Range("N86").Value = -3140

$ws.Range("H89").Value = 217.16667
$ws.Range("I89").Value = 209.4
$ws.Range("J89").Value = 256
$ws.Range("K89").Value = 1884.6
$ws.Range("L89").Value = 2304
$ws.Range("M89").Value = 4043.4
$ws.Range("N89").Value = -14160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8361.125
$ws.Range("I113").Value = 5749.5
$ws.Range("J113").Value = 9231.666999999999
$ws.Range("K113").Value = 5749.5
$ws.Range("L113").Value = 9231.666999999999
$ws.Range("M113").Value = -3579.5
$ws.Range("N113").Value = -13571.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1706
$ws.Range("I22").Value = 1191.5
$ws.Range("J22").Value = 3249.5
$ws.Range("K22").Value = 1191.5
$ws.Range("L22").Value = 3249.5
$ws.Range("M22").Value = -896.5
$ws.Range("N22").Value = -3839.5

$ws.Range("H27").Value = 1706
$ws.Range("I27").Value = 1191.5
$ws.Range("J27").Value = 3249.5
$ws.Range("K27").Value = 1191.5
$ws.Range("L27").Value = 3249.5
$ws.Range("M27").Value = -1084.5
$ws.Range("N27").Value = -3463.5

$ws.Range("H61").Value = 8000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 8000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -8404

$ws.Range("H113").Value = 8000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 8000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -12340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 42273.5
$ws.Range("J52").Value = 83047
$ws.Range("L52").Value = 83047
$ws.Range("N52").Value = -83499

$ws.Range("H113").Value = 704.3158
$ws.Range("I113").Value = 707.0833
$ws.Range("J113").Value = 699.5714
$ws.Range("K113").Value = 2121.2499
$ws.Range("L113").Value = 2098.7142
$ws.Range("M113").Value = 48.7501000000002
$ws.Range("N113").Value = -6438.7142
